# Update loading-percent results for the 380 kV case (Case_5_211).
# Only numeric values in columns B-O, rows 2-25 change; row 1 (headers)
# and columns A, G, I, K, L (all zero) are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 13.70159836729662
$ws.Range("C2").Value = 7.027492222474809
$ws.Range("D2").Value = 9.216955118609004
$ws.Range("E2").Value = 13.49660647338201
$ws.Range("F2").Value = 31.84037344995381
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("J2").Value = 9.917056146039144
$ws.Range("M2").Value = 16.50538171259721
$ws.Range("N2").Value = 17.78526280807522
$ws.Range("O2").Value = 23.84670232787553

# Row 3
$ws.Range("B3").Value = 13.23921161019979
$ws.Range("C3").Value = 6.61474959674558
$ws.Range("D3").Value = 9.207031294696169
$ws.Range("E3").Value = 13.51183840921207
$ws.Range("F3").Value = 31.85815364066604
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("J3").Value = 9.941423444315348
$ws.Range("M3").Value = 16.37104002451241
$ws.Range("N3").Value = 17.84211868669211
$ws.Range("O3").Value = 23.89158391864096

# Row 4
$ws.Range("B4").Value = 12.94898574018271
$ws.Range("C4").Value = 6.347681334806357
$ws.Range("D4").Value = 9.202121101417852
$ws.Range("E4").Value = 13.52329630282902
$ws.Range("F4").Value = 31.87758136147765
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("J4").Value = 9.957505084942525
$ws.Range("M4").Value = 16.29064116373495
$ws.Range("N4").Value = 17.87886909178363
$ws.Range("O4").Value = 23.92519417823688

# Row 5
$ws.Range("B5").Value = 12.82931357957168
$ws.Range("C5").Value = 6.235484225123223
$ws.Range("D5").Value = 9.20041941785751
$ws.Range("E5").Value = 13.52849510792066
$ws.Range("F5").Value = 31.88763557887645
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("J5").Value = 9.964340431528109
$ws.Range("M5").Value = 16.25843208630272
$ws.Range("N5").Value = 17.89430915444407
$ws.Range("O5").Value = 23.94040852223292

# Row 6
$ws.Range("B6").Value = 12.80936313616835
$ws.Range("C6").Value = 6.216652997676183
$ws.Range("D6").Value = 9.200154978009264
$ws.Range("E6").Value = 13.52939035696519
$ws.Range("F6").Value = 31.88943405396449
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("J6").Value = 9.965492473678287
$ws.Range("M6").Value = 16.25311809499796
$ws.Range("N6").Value = 17.896901021272
$ws.Range("O6").Value = 23.94302639897434

# Row 7
$ws.Range("B7").Value = 12.9473772216934
$ws.Range("C7").Value = 6.346181730231794
$ws.Range("D7").Value = 9.20209693799154
$ws.Range("E7").Value = 13.52336427113686
$ws.Range("F7").Value = 31.8777083071037
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("J7").Value = 9.957596126811907
$ws.Range("M7").Value = 16.2902044994212
$ws.Range("N7").Value = 17.87907544176347
$ws.Range("O7").Value = 23.9253932241736

# Row 8
$ws.Range("B8").Value = 13.54358551243868
$ws.Range("C8").Value = 6.88804511493504
$ws.Range("D8").Value = 9.213288936088025
$ws.Range("E8").Value = 13.50142160930454
$ws.Range("F8").Value = 31.84473635418437
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("J8").Value = 9.925225716085011
$ws.Range("M8").Value = 16.45864523580778
$ws.Range("N8").Value = 17.80448536051652
$ws.Range("O8").Value = 23.8609189049714

# Row 9
$ws.Range("B9").Value = 14.65525163213464
$ws.Range("C9").Value = 7.879333994089513
$ws.Range("D9").Value = 9.244541137733833
$ws.Range("E9").Value = 13.47508629247733
$ws.Range("F9").Value = 31.84768058907801
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("J9").Value = 9.870620505705526
$ws.Range("M9").Value = 16.80419742723292
$ws.Range("N9").Value = 17.67276802439003
$ws.Range("O9").Value = 23.78267576154566

# Row 10
$ws.Range("B10").Value = 15.42820276905094
$ws.Range("C10").Value = 8.53974006297636
$ws.Range("D10").Value = 9.273056944259434
$ws.Range("E10").Value = 13.46589525205465
$ws.Range("F10").Value = 31.89108226027658
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("J10").Value = 9.835892442084669
$ws.Range("M10").Value = 17.0655474603299
$ws.Range("N10").Value = 17.5847958208958
$ws.Range("O10").Value = 23.75476397642601

# Row 11
$ws.Range("B11").Value = 15.76878892440211
$ws.Range("C11").Value = 8.822712928729459
$ws.Range("D11").Value = 9.287207774848181
$ws.Range("E11").Value = 13.46391317374797
$ws.Range("F11").Value = 31.91975981593078
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("J11").Value = 9.821260273851792
$ws.Range("M11").Value = 17.18566852561822
$ws.Range("N11").Value = 17.54667092570332
$ws.Range("O11").Value = 23.74851803715001

# Row 12
$ws.Range("B12").Value = 15.89606419454946
$ws.Range("C12").Value = 8.927365604990358
$ws.Range("D12").Value = 9.292733219767184
$ws.Range("E12").Value = 13.46347804219937
$ws.Range("F12").Value = 31.93189976828937
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("J12").Value = 9.815886801823526
$ws.Range("M12").Value = 17.23129868414029
$ws.Range("N12").Value = 17.53250530245571
$ws.Range("O12").Value = 23.7470819861285

# Row 13
$ws.Range("B13").Value = 15.86873030120435
$ws.Range("C13").Value = 8.904937991620347
$ws.Range("D13").Value = 9.291535840641556
$ws.Range("E13").Value = 13.46355773998241
$ws.Range("F13").Value = 31.92922835016706
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("J13").Value = 9.817036633058361
$ws.Range("M13").Value = 17.22146560087803
$ws.Range("N13").Value = 17.53554406464108
$ws.Range("O13").Value = 23.74734992456337

# Row 14
$ws.Range("B14").Value = 15.77929454209841
$ws.Range("C14").Value = 8.831372892446671
$ws.Range("D14").Value = 9.287659028824089
$ws.Range("E14").Value = 13.46387105822343
$ws.Range("F14").Value = 31.92073294588644
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("J14").Value = 9.820814841175848
$ws.Range("M14").Value = 17.18941987769654
$ws.Range("N14").Value = 17.54550007635293
$ws.Range("O14").Value = 23.74838126465533

# Row 15
$ws.Range("B15").Value = 15.72428845358632
$ws.Range("C15").Value = 8.785986380176464
$ws.Range("D15").Value = 9.28530601412052
$ws.Range("E15").Value = 13.46410402904824
$ws.Range("F15").Value = 31.91569585792097
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("J15").Value = 9.823150898063622
$ws.Range("M15").Value = 17.16980854460861
$ws.Range("N15").Value = 17.55163374443551
$ws.Range("O15").Value = 23.74913402581249

# Row 16
$ws.Range("B16").Value = 15.40571287667445
$ws.Range("C16").Value = 8.520896501402122
$ws.Range("D16").Value = 9.272155649352282
$ws.Range("E16").Value = 13.46606897145912
$ws.Range("F16").Value = 31.88938755180861
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("J16").Value = 9.836872128018744
$ws.Range("M16").Value = 17.05771912746847
$ws.Range("N16").Value = 17.58732540726872
$ws.Range("O16").Value = 23.75530211896363

# Row 17
$ws.Range("B17").Value = 15.20736817297315
$ws.Range("C17").Value = 8.353807336579123
$ws.Range("D17").Value = 9.2643883922138
$ws.Range("E17").Value = 13.46783710108349
$ws.Range("F17").Value = 31.87553382904623
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("J17").Value = 9.845588084620179
$ws.Range("M17").Value = 16.98924703479062
$ws.Range("N17").Value = 17.60970554788367
$ws.Range("O17").Value = 23.76073961165903

# Row 18
$ws.Range("B18").Value = 15.09225279293175
$ws.Range("C18").Value = 8.256059224662996
$ws.Range("D18").Value = 9.260031999635238
$ws.Range("E18").Value = 13.46906111426397
$ws.Range("F18").Value = 31.86840678231342
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("J18").Value = 9.850711016601753
$ws.Range("M18").Value = 16.94998167635881
$ws.Range("N18").Value = 17.62275634831692
$ws.Range("O18").Value = 23.76447428955947

# Row 19
$ws.Range("B19").Value = 15.05310291956913
$ws.Range("C19").Value = 8.222681063392095
$ws.Range("D19").Value = 9.258576165841642
$ws.Range("E19").Value = 13.469511123422
$ws.Range("F19").Value = 31.8661382803012
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("J19").Value = 9.85246440904166
$ws.Range("M19").Value = 16.93670840486608
$ws.Range("N19").Value = 17.62720578287464
$ws.Range("O19").Value = 23.76584301547682

# Row 20
$ws.Range("B20").Value = 15.2285900354206
$ws.Range("C20").Value = 8.371764250167773
$ws.Range("D20").Value = 9.265203747641372
$ws.Range("E20").Value = 13.46762745883891
$ws.Range("F20").Value = 31.87692154217528
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("J20").Value = 9.844648899701726
$ws.Range("M20").Value = 16.99652404325765
$ws.Range("N20").Value = 17.60730469339525
$ws.Range("O20").Value = 23.76009792820844

# Row 21
$ws.Range("B21").Value = 15.8056108749562
$ws.Range("C21").Value = 8.853048622011947
$ws.Range("D21").Value = 9.288793236726212
$ws.Range("E21").Value = 13.46377047492574
$ws.Range("F21").Value = 31.92319354395939
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("J21").Value = 9.819700548144983
$ws.Range("M21").Value = 17.19882887701694
$ws.Range("N21").Value = 17.54256839461802
$ws.Range("O21").Value = 23.74805310963628

# Row 22
$ws.Range("B22").Value = 16.17279174590525
$ws.Range("C22").Value = 9.153009879230108
$ws.Range("D22").Value = 9.30518141656262
$ws.Range("E22").Value = 13.46308788232199
$ws.Range("F22").Value = 31.96089505694702
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("J22").Value = 9.804371062511388
$ws.Range("M22").Value = 17.33186400314051
$ws.Range("N22").Value = 17.50184127056336
$ws.Range("O22").Value = 23.74559715702053

# Row 23
$ws.Range("B23").Value = 15.97776276848817
$ws.Range("C23").Value = 8.994247231359957
$ws.Range("D23").Value = 9.296346804695769
$ws.Range("E23").Value = 13.46328428979947
$ws.Range("F23").Value = 31.940092231645
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("J23").Value = 9.81246349521086
$ws.Range("M23").Value = 17.26079715953213
$ws.Range("N23").Value = 17.52343366830398
$ws.Range("O23").Value = 23.74641205732079

# Row 24
$ws.Range("B24").Value = 15.21899900040457
$ws.Range("C24").Value = 8.363651182573459
$ws.Range("D24").Value = 9.264834785404759
$ws.Range("E24").Value = 13.46772159172883
$ws.Range("F24").Value = 31.87629154752505
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("J24").Value = 9.845073156273228
$ws.Range("M24").Value = 16.99323379312544
$ws.Range("N24").Value = 17.60838954597402
$ws.Range("O24").Value = 23.76038613765767

# Row 25
$ws.Range("B25").Value = 14.3616346550579
$ws.Range("C25").Value = 7.620626358039356
$ws.Range("D25").Value = 9.235102014581944
$ws.Range("E25").Value = 13.48042503715314
$ws.Range("F25").Value = 31.83963885888508
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("J25").Value = 9.884444702569386
$ws.Range("M25").Value = 16.70927034934614
$ws.Range("N25").Value = 17.70685048885841
$ws.Range("O25").Value = 23.79866109548798
